# Apply the "금융시장 동향" daily-date roll-forward edit:
#  - 금일보고 sheet: bump the report-date banner from 2021.09.28 to 2021.09.29
#  - History sheet: clear out the (now stale) Hang Seng "J" column daily figures
#  - refresh the selected cell / active sheet bookkeeping for both sheets

$wb = $excel.ActiveWorkbook

$wsToday   = $wb.Worksheets.Item(1)   # 금일보고
$wsHistory = $wb.Worksheets.Item(3)   # History

# 1. Update the title banner (merged G1:N3) with the new reporting date.
#    Changing the text also retires the now-unused old string and appends
#    the new one at the end of the shared-string table automatically.
$wsToday.Range("G1").Value = "금융시장 동향 (2021.09.29.)"

# Re-running AutoFit keeps row 1's height attribute from drifting away from
# the sheet default after the text content changed.
$wsToday.Rows.Item(1).EntireRow.AutoFit()

# 2. Clear the stale daily values in column J ("항셍") of the History sheet
#    for every row that still has a number in it.
$rowsToClear = @(8,10,15,25,30,34,39,44,49,54,59,64,69,74,79,84,89,94,98,103,112,113,115,116,117,118,119,120,121,122,123,124,125,126,127,128,129,130,131,132,133,134,135,136,137,138,139,140,141,142,143,144,145,146,147,148,149,150,151,152,153,154,155,156,157,158,162,163,164,165,166,167,168,169,170,171,172,173,174,175,176,177,178,179,180,181,182,183,184,185,186,187,188)

foreach ($r in $rowsToClear) {
    $wsHistory.Range("J" + $r).ClearContents()
}

# 3. Refresh the saved selection / active-sheet state.
#    History sheet loses focus (its tabSelected flag) while 금일보고 keeps it,
#    so activate History first and 금일보고 last.
$wsHistory.Activate()
$wsHistory.Range("J8:J328").Select()

$wsToday.Activate()
$wsToday.Range("A13").Select()
